$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "El peso del corazón"
$ws.Range("B11").Value = "Rosa Montero"
$ws.Range("C11").Value = ""
